$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Task"
$ws.Range("B1").Value = "Expected"
$ws.Range("C1").Value = "Actual"
$ws.Range("D1").Value = "Difference"

# Task rows: Name, Expected, Actual
$tasks = @(
    @("Documentation/SA demo", 24, 24),
    @("Github setup, Etc.", 6, 18),
    @("Write Public functions, basic structure", 6, 7),
    @("Make inert obstacles", 6, 4),
    @("Make Simple Interactables", 3, 4),
    @("Make Level Exit/Dodo Egg", 3, 6),
    @("Make Optional Obstacles", 6, 4),
    @("Integrate other's parts that were not previously working", 3, 7),
    @("Final testing and build", 3, 3)
)

$row = 2
foreach ($task in $tasks) {
    $ws.Cells.Item($row, 1).Value = $task[0]
    $ws.Cells.Item($row, 2).Value = $task[1]
    $ws.Cells.Item($row, 3).Value = $task[2]
    $ws.Cells.Item($row, 4).Formula = "=B$row-C$row"
    $row = $row + 1
}

# Totals row
$ws.Range("A11").Value = "Total"
$ws.Range("B11").Formula = "=SUM(B2:B10)"
$ws.Range("C11").Formula = "=SUM(C2:C10)"
$ws.Range("D11").Formula = "=SUM(D2:D10)"

# Column width / selection cosmetics
$ws.Columns.Item(1).ColumnWidth = 51.7109375
$null = $ws.Range("E17").Select()

# Page was set up for printing (portrait) as part of this edit
$ws.PageSetup.Orientation = 1
